$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells should look like the existing header row (bold, bordered, centered),
# so copy formatting from an existing header cell (AC1) onto AD1:AF1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Add header labels for the new "team record" columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the (constant) team record for every player row (2-49)
for ($r = 2; $r -le 49; $r++) {
    $ws.Cells.Item($r, 30).Value = 94  # AD: Wins
    $ws.Cells.Item($r, 31).Value = 68  # AE: Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF: Ties
}
